$wb = $excel.ActiveWorkbook

# ALC row 98
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H98").Value = 631.25
$ws.Range("I98").Value = 631.25
$ws.Range("J98").Value = 0
$ws.Range("K98").Value = 631.25
$ws.Range("L98").Value = 0
$ws.Range("M98").Value = 866.75
$ws.Range("N98").ClearContents()

# ALC row 112
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H112").Value = 1491.9166
$ws.Range("J112").Value = 1655.8889
$ws.Range("L112").Value = 4967.6667
$ws.Range("N112").Value = -7183.6667

# ALC row 116
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H116").Value = 5557640
$ws.Range("I116").Value = 6912567
$ws.Range("K116").Value = 6912567
$ws.Range("M116").Value = -6909125

# ALC row 122
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H122").Value = 631.25
$ws.Range("I122").Value = 631.25
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 1893.75
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = 556.25
$ws.Range("N122").ClearContents()

# ALC row 127
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H127").Value = 981.8333
$ws.Range("I127").Value = 724.6667
$ws.Range("K127").Value = 2174.0001
$ws.Range("M127").Value = 2785.9999

# ALC row 129
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H129").Value = 844
$ws.Range("J129").Value = 1106.5714
$ws.Range("L129").Value = 3319.7142
$ws.Range("N129").Value = -13319.7142

# ALC row 132
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H132").Value = 3775.6667
$ws.Range("I132").Value = 3705.9792
$ws.Range("J132").Value = 4333.1665
$ws.Range("K132").Value = 11117.9376
$ws.Range("L132").Value = 12999.4995
$ws.Range("M132").Value = -8587.937600000001
$ws.Range("N132").Value = -18059.4995

# ALC row 138
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 1703.096
$ws.Range("I138").Value = 560.6
$ws.Range("J138").Value = 3539.25
$ws.Range("K138").Value = 1681.8
$ws.Range("L138").Value = 10617.75
$ws.Range("M138").Value = 3458.2
$ws.Range("N138").Value = -20897.75

# ALC row 141
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H141").Value = 3400.1404
$ws.Range("I141").Value = 1512.2291
$ws.Range("J141").Value = 13469
$ws.Range("K141").Value = 4536.6873
$ws.Range("L141").Value = 40407
$ws.Range("M141").Value = 643.3127000000004
$ws.Range("N141").Value = -50767

# ARM row 2
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 778.8570999999999
$ws.Range("I2").Value = 766.5
$ws.Range("J2").Value = 853
$ws.Range("K2").Value = 766.5
$ws.Range("L2").Value = 853
$ws.Range("M2").Value = -653.5
$ws.Range("N2").Value = -1079

# ARM row 45
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 1743.3684
$ws.Range("I45").Value = 1744.9333
$ws.Range("J45").Value = 1737.5
$ws.Range("K45").Value = 1744.9333
$ws.Range("L45").Value = 1737.5
$ws.Range("M45").Value = -1367.9333
$ws.Range("N45").Value = -2491.5

# ARM row 74
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 9262583
$ws.Range("I74").Value = 12502921
$ws.Range("J74").Value = 4475.2856
$ws.Range("K74").Value = 12502921
$ws.Range("L74").Value = 4475.2856
$ws.Range("M74").Value = -12502047
$ws.Range("N74").Value = -6223.2856

# ARM row 77
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H77").Value = 9262583
$ws.Range("I77").Value = 12502921
$ws.Range("J77").Value = 4475.2856
$ws.Range("K77").Value = 62514605
$ws.Range("L77").Value = 22376.428
$ws.Range("M77").Value = -62510237
$ws.Range("N77").Value = -31112.428

# ARM row 116
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H116").Value = 778.8570999999999
$ws.Range("I116").Value = 766.5
$ws.Range("J116").Value = 853
$ws.Range("K116").Value = 766.5
$ws.Range("L116").Value = 853
$ws.Range("M116").Value = 1527.5
$ws.Range("N116").Value = -5441

# ARM row 132
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 1647.9556
$ws.Range("I132").Value = 1492.5186
$ws.Range("J132").Value = 1881.1111
$ws.Range("K132").Value = 4477.5558
$ws.Range("L132").Value = 5643.3333
$ws.Range("M132").Value = -1947.5558
$ws.Range("N132").Value = -10703.3333

# BSM row 3
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 778.8570999999999
$ws.Range("I3").Value = 766.5
$ws.Range("J3").Value = 853
$ws.Range("K3").Value = 766.5
$ws.Range("L3").Value = 853
$ws.Range("M3").Value = -652.5
$ws.Range("N3").Value = -1081

# BSM row 26
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H26").Value = 14490.333
$ws.Range("I26").Value = 10235.5
$ws.Range("J26").Value = 23000
$ws.Range("K26").Value = 10235.5
$ws.Range("L26").Value = 23000
$ws.Range("M26").Value = -9943.5
$ws.Range("N26").Value = -23584

# CRP row 43
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H43").Value = 80000
$ws.Range("J43").Value = 80000
$ws.Range("L43").Value = 80000
$ws.Range("N43").Value = -80368

# CRP row 101
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H101").Value = 80000
$ws.Range("J101").Value = 80000
$ws.Range("L101").Value = 80000
$ws.Range("N101").Value = -86490

# CRP row 134
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H134").Value = 2138.5952
$ws.Range("I134").Value = 2409.4138
$ws.Range("J134").Value = 1534.4615
$ws.Range("K134").Value = 7228.241399999999
$ws.Range("L134").Value = 4603.3845
$ws.Range("M134").Value = -4693.241399999999
$ws.Range("N134").Value = -9673.3845

# CUL row 122
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H122").Value = 488.6316
$ws.Range("I122").Value = 324.91666
$ws.Range("J122").Value = 769.2857
$ws.Range("K122").Value = 2924.24994
$ws.Range("L122").Value = 6923.571300000001
$ws.Range("M122").Value = -474.2499399999997
$ws.Range("N122").Value = -11823.5713

# CUL row 131
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 1026.9445
$ws.Range("J131").Value = 1067.5454
$ws.Range("L131").Value = 3202.6362
$ws.Range("N131").Value = -13282.6362

# CUL row 132
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H132").Value = 812.1429000000001
$ws.Range("I132").Value = 830
$ws.Range("J132").Value = 780
$ws.Range("K132").Value = 7470
$ws.Range("L132").Value = 7020
$ws.Range("M132").Value = -4940
$ws.Range("N132").Value = -12080

# CUL row 137
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H137").Value = 2890.44
$ws.Range("J137").Value = 5635.5454
$ws.Range("L137").Value = 16906.6362
$ws.Range("N137").Value = -27106.6362

# LTW row 132
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 11369127
$ws.Range("I132").Value = 36778084
$ws.Range("J132").Value = 1960.9474
$ws.Range("K132").Value = 110334252
$ws.Range("L132").Value = 5882.8422
$ws.Range("M132").Value = -110331722
$ws.Range("N132").Value = -10942.8422

# LTW row 136
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 3019.0164
$ws.Range("I136").Value = 3541.932
$ws.Range("J136").Value = 1665.5883
$ws.Range("K136").Value = 10625.796
$ws.Range("L136").Value = 4996.7649
$ws.Range("M136").Value = -8075.795999999998
$ws.Range("N136").Value = -10096.7649

# WVR row 126
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 1938.359
$ws.Range("I126").Value = 1958.7059
$ws.Range("K126").Value = 5876.1177
$ws.Range("M126").Value = -3406.1177
